# Upload new version with timestamp
# Adds two new shortage-report rows (GLIMET FORTE, and the original HEPTA
# item moved down), re-labels the former row 7 item as CARVID, refreshes
# the running total and bumps the generated-at timestamp by one minute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: insert two blank rows at row 8 (old row 8 -> 10, old row
#    9 -> 11). Excel shifts existing merged ranges along automatically.
# ---------------------------------------------------------------------
$ws.Range("A8:A9").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2) Clone row 7's cell formatting onto the two freshly inserted rows so
#    the new item rows look just like the existing item row.
# ---------------------------------------------------------------------
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)

# Row heights for the two new item rows (row 8 keeps the old total row's
# height slot, row 9 matches the template row's height).
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5

# ---------------------------------------------------------------------
# 3) Re-merge the split cells for the two new rows (PasteSpecial of
#    formats does not recreate merged ranges).
# ---------------------------------------------------------------------
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# ---------------------------------------------------------------------
# Helper-style writes: several columns carry a numeric-format style but
# must hold literal text (e.g. "45.0000"), matching the source report's
# behaviour. Flip the cell to text, write the value, then restore the
# original number format so the style id used on disk is unchanged.
# ---------------------------------------------------------------------
function Set-TextValue($rng, $value) {
    $fmt = $rng.NumberFormat()
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 4) Row 7 now describes item 1 = CARVID 6.25MG 30TAB.
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("C7") "CARVID 6.25MG 30TAB"
Set-TextValue $ws.Range("H7") "0:1"
Set-TextValue $ws.Range("L7") "0"
Set-TextValue $ws.Range("N7") "45.00"
Set-TextValue $ws.Range("P7") "45.0000"
Set-TextValue $ws.Range("Q7") "1:0"

# ---------------------------------------------------------------------
# 5) Row 8 is the new item 2 = GLIMET FORTE 5/800 MG 30 F.C.TAB.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 2
Set-TextValue $ws.Range("C8") "GLIMET FORTE 5/800 MG 30 F.C.TAB."
Set-TextValue $ws.Range("H8") "0:1"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "66.00"
Set-TextValue $ws.Range("P8") "66.0000"
Set-TextValue $ws.Range("Q8") "1:0"

# ---------------------------------------------------------------------
# 6) Row 9 is the new item 3 = HEPTA PANTHENOL HAIR CREAM 100 GM (the
#    item that used to be row 7, now pushed to the bottom of the list).
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 3
Set-TextValue $ws.Range("C9") "HEPTA PANTHENOL HAIR CREAM 100 GM"
Set-TextValue $ws.Range("H9") "0:0"
Set-TextValue $ws.Range("L9") "1"
Set-TextValue $ws.Range("N9") "149.00"
Set-TextValue $ws.Range("P9") "149.0000"
Set-TextValue $ws.Range("Q9") "1:0"

# ---------------------------------------------------------------------
# 7) Row 10 (formerly row 8) holds the grand total of the three prices.
# ---------------------------------------------------------------------
$ws.Range("P10").Value = 260

# ---------------------------------------------------------------------
# 8) Row 11 (formerly row 9) is the footer; only the generated-at
#    timestamp changes, from 9:34 AM to 9:35 AM.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Sunday, 7 September, 2025 9:35 AM"
